# Auto-generated edit script applying numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H-N)
# across multiple job sheets, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3973.25
$ws.Range("I19").Value = 4980.6665
$ws.Range("K19").Value = 4980.6665
$ws.Range("M19").Value = -4805.6665
$ws.Range("H34").Value = 2230
$ws.Range("I34").Value = 876
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 876
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -673
$ws.Range("N34").Value = -9406
$ws.Range("H36").Value = 2230
$ws.Range("I36").Value = 876
$ws.Range("J36").Value = 9000
$ws.Range("K36").Value = 876
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = -161
$ws.Range("N36").Value = -10430
$ws.Range("H97").Value = 1849.5
$ws.Range("I97").Value = 1699
$ws.Range("K97").Value = 5097
$ws.Range("M97").Value = -4601
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 7080.727
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 7861
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 7861
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -14369
$ws.Range("H138").Value = 3534.0889
$ws.Range("I138").Value = 2533.9546
$ws.Range("J138").Value = 4490.7393
$ws.Range("K138").Value = 7601.8638
$ws.Range("L138").Value = 13472.2179
$ws.Range("M138").Value = -2461.8638
$ws.Range("N138").Value = -23752.2179

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2504.9333
$ws.Range("I102").Value = 2504.9333
$ws.Range("K102").Value = 2504.9333
$ws.Range("M102").Value = -882.9333000000001
$ws.Range("H122").Value = 1600.9584
$ws.Range("J122").Value = 2775.75
$ws.Range("L122").Value = 8327.25
$ws.Range("N122").Value = -13227.25
$ws.Range("H132").Value = 2307.6558
$ws.Range("I132").Value = 1874.1818
$ws.Range("K132").Value = 5622.5454
$ws.Range("M132").Value = -3092.5454
$ws.Range("H135").Value = 65228.4
$ws.Range("J135").Value = 65228.4
$ws.Range("L135").Value = 65228.4
$ws.Range("N135").Value = -75368.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1555.5555
$ws.Range("I99").Value = 1492.1666
$ws.Range("K99").Value = 1492.1666
$ws.Range("M99").Value = 5.833399999999983

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27963.285
$ws.Range("I31").Value = 2880.375
$ws.Range("K31").Value = 2880.375
$ws.Range("M31").Value = -2585.375
$ws.Range("H34").Value = 27963.285
$ws.Range("I34").Value = 2880.375
$ws.Range("K34").Value = 2880.375
$ws.Range("M34").Value = -2678.375
$ws.Range("H111").Value = 98583.336
$ws.Range("J111").Value = 98583.336
$ws.Range("L111").Value = 98583.336
$ws.Range("N111").Value = -106763.336
$ws.Range("H132").Value = 5707.7144
$ws.Range("I132").Value = 512
$ws.Range("J132").Value = 6573.6665
$ws.Range("K132").Value = 1536
$ws.Range("L132").Value = 19720.9995
$ws.Range("M132").Value = 994
$ws.Range("N132").Value = -24780.9995
$ws.Range("H141").Value = 174369.5
$ws.Range("J141").Value = 174369.5
$ws.Range("L141").Value = 174369.5
$ws.Range("N141").Value = -184729.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 394.33334
$ws.Range("I23").Value = 338
$ws.Range("J23").Value = 405.6
$ws.Range("K23").Value = 1014
$ws.Range("L23").Value = 1216.8
$ws.Range("M23").Value = -779
$ws.Range("N23").Value = -1686.8
$ws.Range("H69").Value = 7503.5
$ws.Range("J69").Value = 9995
$ws.Range("L69").Value = 29985
$ws.Range("N69").Value = -31607
$ws.Range("H72").Value = 7503.5
$ws.Range("J72").Value = 9995
$ws.Range("L72").Value = 89955
$ws.Range("N72").Value = -98067
$ws.Range("H75").Value = 83339010
$ws.Range("J75").Value = 8274.625
$ws.Range("L75").Value = 24823.875
$ws.Range("N75").Value = -26819.875
$ws.Range("H78").Value = 83339010
$ws.Range("J78").Value = 8274.625
$ws.Range("L78").Value = 74471.625
$ws.Range("N78").Value = -84455.625
$ws.Range("H86").Value = 1997.875
$ws.Range("I86").Value = 556.4
$ws.Range("J86").Value = 4400.3335
$ws.Range("K86").Value = 1669.2
$ws.Range("L86").Value = 13201.0005
$ws.Range("M86").Value = -483.1999999999998
$ws.Range("N86").Value = -15573.0005
$ws.Range("H89").Value = 1997.875
$ws.Range("I89").Value = 556.4
$ws.Range("J89").Value = 4400.3335
$ws.Range("K89").Value = 5007.599999999999
$ws.Range("L89").Value = 39603.0015
$ws.Range("M89").Value = 920.4000000000005
$ws.Range("N89").Value = -51459.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 164999.5
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H97").Value = 1797.75
$ws.Range("J97").Value = 2442.8
$ws.Range("L97").Value = 2442.8
$ws.Range("N97").Value = -3434.8
$ws.Range("H102").Value = 2814
$ws.Range("I102").Value = 2814
$ws.Range("K102").Value = 2814
$ws.Range("M102").Value = -1192
$ws.Range("H107").Value = 1452.2222
$ws.Range("I107").Value = 760.44446
$ws.Range("K107").Value = 760.44446
$ws.Range("M107").Value = 1159.55554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3112.96
$ws.Range("I82").Value = 2872.5293
$ws.Range("J82").Value = 3623.875
$ws.Range("K82").Value = 2872.5293
$ws.Range("L82").Value = 3623.875
$ws.Range("M82").Value = -2511.5293
$ws.Range("N82").Value = -4345.875
$ws.Range("H85").Value = 3112.96
$ws.Range("I85").Value = 2872.5293
$ws.Range("J85").Value = 3623.875
$ws.Range("K85").Value = 2872.5293
$ws.Range("L85").Value = 3623.875
$ws.Range("M85").Value = -1624.5293
$ws.Range("N85").Value = -6119.875
$ws.Range("H132").Value = 2630.9062
$ws.Range("I132").Value = 1277.56
$ws.Range("K132").Value = 3832.68
$ws.Range("M132").Value = -1302.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 9000
$ws.Range("J25").Value = 9000
$ws.Range("L25").Value = 9000
$ws.Range("N25").Value = -9586
$ws.Range("H132").Value = 6677.649
$ws.Range("I132").Value = 4711.694
$ws.Range("J132").Value = 18719.125
$ws.Range("K132").Value = 14135.082
$ws.Range("L132").Value = 56157.375
$ws.Range("M132").Value = -11605.082
$ws.Range("N132").Value = -61217.375
